# Add a flexible "Mark" column to the generic parser template.
# This inserts a new column before column I (Stock) on the "Template"
# sheet, pushing every column from I onward one slot to the right, and
# labels the new column header "Mark".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Insert a new blank column at I - everything currently in/after I
# (Stock, From Tank, To Tank, Tagger, ...) shifts right by one.
$ws.Columns("I").Insert()

# New column inherits column H's formatting but, like Excel's own
# insert-column behaviour, starts out without an explicit "best fit"
# width flag since it has no content yet to size to - give it an
# explicit width close to its neighbour's.
$ws.Columns("I").ColumnWidth = 8

# Label the new column.
$ws.Range("I1").Value = "Mark"

# Leave the new header cell selected, matching the saved selection.
$ws.Range("I1").Select() | Out-Null
